# Applies the "Aldabón-Gemini" rebrand + refreshed study data described in the
# commit "feat: Add Technical Guide, PDF Logic, and Study Data for Cloud".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Comparativa_Final": rename the model column + refresh its figures.
# ---------------------------------------------------------------------------
$wsComp = $wb.Worksheets.Item("Comparativa_Final")

$wsComp.Range("C1").Value = "Aldabón-Gemini"

# Partido | CIS (Oficial) | Aldabón-Gemini | Diferencia
$wsComp.Range("B2").Value = 23
$wsComp.Range("C2").Value = 33.8
$wsComp.Range("D2").Value = 10.8

$wsComp.Range("B3").Value = 31.7
$wsComp.Range("C3").Value = 27.5
$wsComp.Range("D3").Value = -4.2

$wsComp.Range("B4").Value = 17.7
$wsComp.Range("C4").Value = 18.7
$wsComp.Range("D4").Value = 1

$wsComp.Range("B5").Value = 7.2
$wsComp.Range("C5").Value = 6.4
$wsComp.Range("D5").Value = -0.8

# ---------------------------------------------------------------------------
# Sheet "Detalle_Calculos": refresh the underlying calculation columns.
# ---------------------------------------------------------------------------
$wsDet = $wb.Worksheets.Item("Detalle_Calculos")

# Recuerdo_CIS | K_Ponderacion | Voto_Simpatia_CIS | Ajuste_Fidelidad | Final_%
$wsDet.Range("C2").Value = 675.9490378237615
$wsDet.Range("D2").Value = 1.378
$wsDet.Range("E2").Value = 15.9
$wsDet.Range("G2").Value = 33.8

$wsDet.Range("C3").Value = 1183.719221212038
$wsDet.Range("D3").Value = 0.753
$wsDet.Range("E3").Value = 23.2
$wsDet.Range("G3").Value = 27.5

$wsDet.Range("C4").Value = 350.4993872341499
$wsDet.Range("D4").Value = 0.995
$wsDet.Range("E4").Value = 13.8
$wsDet.Range("G4").Value = 18.7

$wsDet.Range("C5").Value = 366.8012529169699
$wsDet.Range("D5").Value = 0.9429999999999999
$wsDet.Range("E5").Value = 4.8
$wsDet.Range("G5").Value = 6.4

$wsDet.Range("C6").Value = 57.28611840005001
$wsDet.Range("D6").Value = 0.9330000000000001
$wsDet.Range("E6").Value = 1.3
$wsDet.Range("G6").Value = 2.1

$wsDet.Range("C7").Value = 46.01147293533
$wsDet.Range("D7").Value = 0.978
$wsDet.Range("E7").Value = 0.7
$wsDet.Range("G7").Value = 1.2

$wsDet.Range("C8").Value = 38.99718611139
$wsDet.Range("D8").Value = 1.01
$wsDet.Range("E8").Value = 1
$wsDet.Range("G8").Value = 1.8

$wsDet.Range("C9").Value = 24.34094112978999
$wsDet.Range("D9").Value = 1.272
$wsDet.Range("E9").Value = 0.5
$wsDet.Range("G9").Value = 1.1

$wsDet.Range("C10").Value = 27.79681177289
$wsDet.Range("D10").Value = 0.607
$wsDet.Range("E10").Value = 0.8
$wsDet.Range("G10").Value = 0.8

$wsDet.Range("C11").Value = 17.0494443116
$wsDet.Range("D11").Value = 0.825
$wsDet.Range("G11").Value = 0.3

$wsDet.Range("C12").Value = 1.0960460554
$wsDet.Range("D12").Value = 5.134
$wsDet.Range("E12").Value = 0.02
$wsDet.Range("G12").Value = 0.2

$wsDet.Range("C13").Value = 24.0848633404
$wsDet.Range("D13").Value = 0.8179999999999999
